{"js": "// Insert three new paragraphs right after the \"Write Up\" title paragraph\n// and before the pre-existing blank paragraph:\n//   1. Intro paragraph about preventing the hero from going through walls.\n//   2. \"So, if this sounds interesting...\" paragraph.\n//   3. \"7 Wall Collisions\" heading (Heading1 style).\n//\n// The inserts are anchored \"Before\" the pre-existing blank (un-styled)\n// paragraph that originally follows \"Write Up\" - inserting there (instead\n// of \"After\" the Title paragraph) keeps the two new body paragraphs free\n// of any inherited \"Title\" formatting, matching the target markup (no\n// <w:pPr> on those two <w:p> elements).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/styleBuiltIn\");\nawait context.sync();\n\n// Find the Title paragraph (\"Write Up\"); the paragraph right after it is\n// the pre-existing blank paragraph we anchor our inserts on.\nconst items = paragraphs.items;\nlet titleIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].styleBuiltIn === \"Title\") {\n    titleIndex = i;\n    break;\n  }\n}\n\nconst blankParagraph = items[titleIndex + 1];\n\nconst introParagraph = blankParagraph.insertParagraph(\n  \"This week, we will be taking a look at how we can prevent the hero object from going through the wall. This will be done through the process of creating both a controller object, and a script. A script is code that can be use globally through out your project, and works very well, when it comes to preventing our little guy from becoming a ghost, and just walking through all the walls that we put up.\",\n  \"Before\"\n);\n\nconst joinParagraph = blankParagraph.insertParagraph(\n  \"So, if this sounds interesting to you then why don\\u2019t you join us for our brand-new article entitled:\",\n  \"Before\"\n);\n\nconst headingParagraph = blankParagraph.insertParagraph(\"7 Wall Collisions\", \"Before\");\nheadingParagraph.style = \"Heading1\";\n\nawait context.sync();\n", "ps1": "# Insert three new paragraphs right after the \"Write Up\" title paragraph\n# and before the pre-existing blank paragraph:\n#   1. Intro paragraph about preventing the hero from going through walls.\n#   2. \"So, if this sounds interesting...\" paragraph.\n#   3. \"7 Wall Collisions\" heading (Heading1 style).\n#\n# The inserts are anchored on the Range of the pre-existing blank (Normal\n# style) paragraph that originally follows \"Write Up\" - inserting there\n# (instead of \"after\" the Title paragraph) keeps the two new body\n# paragraphs free of any inherited Title formatting, matching the target\n# markup (no <w:pPr> on those two <w:p> elements).\n\n$d = $word.ActiveDocument\n\n# Locate the \"Write Up\" Title paragraph, then the blank paragraph right\n# after it (both present in the original document).\n$titleIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  if ($d.Paragraphs.Item($i).Style.NameLocal -eq \"Title\") {\n    $titleIndex = $i\n    break\n  }\n}\n\n$blankIndex = $titleIndex + 1\n$blankParagraph = $d.Paragraphs.Item($blankIndex)\n$anchorRange = $blankParagraph.Range\n\n# Insert three empty paragraphs immediately before the blank paragraph,\n# in document order.\n$anchorRange.InsertParagraphBefore()\n$anchorRange.InsertParagraphBefore()\n$anchorRange.InsertParagraphBefore()\n\n$introParagraph = $d.Paragraphs.Item($blankIndex)\n$introParagraph.Range.Text = \"This week, we will be taking a look at how we can prevent the hero object from going through the wall. This will be done through the process of creating both a controller object, and a script. A script is code that can be use globally through out your project, and works very well, when it comes to preventing our little guy from becoming a ghost, and just walking through all the walls that we put up.\"\n\n$joinParagraph = $d.Paragraphs.Item($blankIndex + 1)\n$joinParagraph.Range.Text = \"So, if this sounds interesting to you then why don\" + [char]8217 + \"t you join us for our brand-new article entitled:\"\n\n$headingParagraph = $d.Paragraphs.Item($blankIndex + 2)\n$headingParagraph.Range.Text = \"7 Wall Collisions\"\n$headingParagraph.Style = \"Heading1\"\n"}
